$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from diff: price/volume refresh + two row
# position swaps (EnergySwap<->Maker at 40/41, OKB<->InjectiveProtocol at 43/44).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.523.47"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.414.90"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.53"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.47"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +6.98%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.421.36"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.006.68"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.98"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.558.33"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.419.37"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.37"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.98"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.50"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +8.42%  "
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.49"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +3.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.13"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.24"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +6.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.44"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.99"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +6.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0763"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.886.68"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.68"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.45"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +9.81%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.74"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0316"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.769"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "322.34"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +8.25%  "
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.855"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.06%  "
